$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs / Pf4 / Procr / ECs ---
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pf4"
$ws.Range("C2").Value = "Procr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 18.359437
$ws.Range("H2").Value = 55.078311
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 7.227365333333334
$ws.Range("N2").Value = 21.682096
$ws.Range("O2").Value = 0.1241091393606781
$ws.Range("P2").Value = 0.1241091393606781
$ws.Range("Q2").Value = 132.6903585133173
$ws.Range("R2").Value = 1194.213226619856
$ws.Range("S2").Value = 0.1241091393606781
$ws.Range("T2").Value = 0.1241091393606781

# --- Row 3: ECs / Pf4 / Procr / FAPs ---
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pf4"
$ws.Range("C3").Value = "Procr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 18.359437
$ws.Range("H3").Value = 55.078311
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 30.598211
$ws.Range("N3").Value = 91.794633
$ws.Range("O3").Value = 0.5254359587541398
$ws.Range("P3").Value = 0.5254359587541398
$ws.Range("Q3").Value = 561.765927167207
$ws.Range("R3").Value = 5055.893344504863
$ws.Range("S3").Value = 0.5254359587541398
$ws.Range("T3").Value = 0.5254359587541398

# --- Row 4 (new row): ECs / Pf4 / Procr / sCs ---
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pf4"
$ws.Range("C4").Value = "Procr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 18.359437
$ws.Range("H4").Value = 55.078311
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 20.408373
$ws.Range("N4").Value = 61.22511899999999
$ws.Range("O4").Value = 0.3504549018851821
$ws.Range("P4").Value = 0.3504549018851821
$ws.Range("Q4").Value = 374.686238366001
$ws.Range("R4").Value = 3372.176145294009
$ws.Range("S4").Value = 0.3504549018851821
$ws.Range("T4").Value = 0.3504549018851821
